$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: remove the standalone "Meta description: ..." paragraph that used
# to sit right under the H1 title.
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Meta description")) {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Step 2: locate the trailing "Please create a cartoon style feature image
# ..." image-prompt paragraph (last paragraph of the body) and, just before
# it, insert a new bold paragraph repeating the page title.
# ---------------------------------------------------------------------------
$promptText = "Please create a cartoon style feature image for the game " + [char]34 + "Da Vinci Ways" + [char]34 + ". The image should feature a happy Maya warrior with glasses. The warrior should be depicted standing in front of a golden painting that contains the reels of the game, with iconic works of Leonardo da Vinci visible in the painting. The warrior should have a joyful expression on their face, holding up a mobile device with the game on it as if they just won a big payout. The image should convey a sense of fun and excitement while also highlighting the game" + [char]39 + "s theme of art and culture."

$promptIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Please create a cartoon style feature image")) {
        $promptIndex = $i
        break
    }
}

if ($promptIndex -ge 2) {
    $titleText = "Play Da Vinci Ways Slot for Free | Review by Experts"

    $anchor = $d.Paragraphs.Item($promptIndex - 1)
    $anchor.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($promptIndex)
    $newPara.Style = "Normal"
    $newPara.Range.Text = $titleText

    $titleRange = $d.Range($newPara.Range.Start, $newPara.Range.Start + $titleText.Length)
    $titleRange.Font.Bold = 1

    $leadIn = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $emptyRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $leadIn.InsertXML($emptyRunXml)
}

# ---------------------------------------------------------------------------
# Step 3: swap the image-prompt copy for the new meta-description style text
# in that final paragraph (keeps its italic run formatting).
# ---------------------------------------------------------------------------
$replaced = $d.Content.Find.Execute(
    $promptText,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our unbiased review of Da Vinci Ways slot game. Play for free and discover its features, symbols, payouts, and mobile compatibility.",
    2)
